# Auto-generated Excel COM-interop script
# Applies numeric corrections to the Leve profit-tracking sheets
# (per-sheet Table_<CODE> data) to match the scheduled-runner update.

$wb = $excel.ActiveWorkbook

# ---- Sheet: ALC ----
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H6").Value = 452.5
$ws.Range("I6").Value = 379.25
$ws.Range("J6").Value = 599
$ws.Range("K6").Value = 1137.75
$ws.Range("L6").Value = 1797
$ws.Range("M6").Value = -1025.75
$ws.Range("N6").Value = -2021
$ws.Range("H17").Value = 950
$ws.Range("J17").Value = 998.3333
$ws.Range("L17").Value = 2994.9999
$ws.Range("N17").Value = -3330.9999
$ws.Range("H69").Value = 42401.3
$ws.Range("J69").Value = 23750
$ws.Range("L69").Value = 71250
$ws.Range("N69").Value = -72998
$ws.Range("H70").Value = 1691.1666
$ws.Range("I70").Value = 1549.6666
$ws.Range("K70").Value = 4648.9998
$ws.Range("M70").Value = -4378.9998
$ws.Range("H72").Value = 42401.3
$ws.Range("J72").Value = 23750
$ws.Range("L72").Value = 213750
$ws.Range("N72").Value = -222486
$ws.Range("H73").Value = 1691.1666
$ws.Range("I73").Value = 1549.6666
$ws.Range("K73").Value = 4648.9998
$ws.Range("M73").Value = -3712.9998
$ws.Range("H80").Value = 1308
$ws.Range("I80").Value = 346.66666
$ws.Range("J80").Value = 2750
$ws.Range("K80").Value = 1039.99998
$ws.Range("L80").Value = 8250
$ws.Range("M80").Value = -41.99998000000005
$ws.Range("N80").Value = -10246
$ws.Range("H83").Value = 1308
$ws.Range("I83").Value = 346.66666
$ws.Range("J83").Value = 2750
$ws.Range("K83").Value = 3119.99994
$ws.Range("L83").Value = 24750
$ws.Range("M83").Value = 1872.00006
$ws.Range("N83").Value = -34734
$ws.Range("H132").Value = 2737.818
$ws.Range("I132").Value = 2876.6
$ws.Range("K132").Value = 8629.799999999999
$ws.Range("M132").Value = -6099.799999999999

# ---- Sheet: ARM ----
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H61").Value = 2578
$ws.Range("I61").Value = 2756
$ws.Range("K61").Value = 2756
$ws.Range("M61").Value = -2544
$ws.Range("H63").Value = 14112.875
$ws.Range("I63").Value = 13817.333
$ws.Range("J63").Value = 14999.5
$ws.Range("K63").Value = 13817.333
$ws.Range("L63").Value = 14999.5
$ws.Range("M63").Value = -13131.333
$ws.Range("N63").Value = -16371.5
$ws.Range("H64").Value = 0
$ws.Range("J64").Value = 0
$ws.Range("L64").ClearContents()
$ws.Range("N64").Value = 0
$ws.Range("H66").Value = 14112.875
$ws.Range("I66").Value = 13817.333
$ws.Range("J66").Value = 14999.5
$ws.Range("K66").Value = 69086.66500000001
$ws.Range("L66").Value = 74997.5
$ws.Range("M66").Value = -65654.66500000001
$ws.Range("N66").Value = -81861.5
$ws.Range("H67").Value = 0
$ws.Range("J67").Value = 0
$ws.Range("L67").ClearContents()
$ws.Range("N67").Value = 0
$ws.Range("H74").Value = 1581.5
$ws.Range("I74").Value = 1297.8
$ws.Range("K74").Value = 1297.8
$ws.Range("M74").Value = -423.8
$ws.Range("H77").Value = 1581.5
$ws.Range("I77").Value = 1297.8
$ws.Range("K77").Value = 6489
$ws.Range("M77").Value = -2121
$ws.Range("H132").Value = 1233.5
$ws.Range("I132").Value = 1200.2
$ws.Range("K132").Value = 3600.6
$ws.Range("M132").Value = -1070.6
$ws.Range("H136").Value = 2578
$ws.Range("I136").Value = 2756
$ws.Range("K136").Value = 8268
$ws.Range("M136").Value = -5718
$ws.Range("H139").Value = 0
$ws.Range("J139").Value = 0
$ws.Range("L139").ClearContents()
$ws.Range("N139").Value = 0

# ---- Sheet: BSM ----
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H107").Value = 1791.8889
$ws.Range("I107").Value = 1518.2858
$ws.Range("J107").Value = 2749.5
$ws.Range("K107").Value = 1518.2858
$ws.Range("L107").Value = 2749.5
$ws.Range("M107").Value = 401.7141999999999
$ws.Range("N107").Value = -6589.5
$ws.Range("H134").Value = 2400.6
$ws.Range("I134").Value = 2400.6
$ws.Range("J134").Value = 0
$ws.Range("K134").Value = 7201.799999999999
$ws.Range("L134").Value = 0
$ws.Range("M134").ClearContents()
$ws.Range("N134").Value = -4666.799999999999
$ws.Range("H9").Value = 90000
$ws.Range("J9").Value = 90000
$ws.Range("L9").Value = 90000
$ws.Range("N9").Value = -90336

# ---- Sheet: CRP ----
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H58").Value = 1663.8
$ws.Range("J58").Value = 2164.5
$ws.Range("L58").Value = 2164.5
$ws.Range("N58").Value = -2570.5
$ws.Range("H94").Value = 104353.73
$ws.Range("I94").Value = 223564.8
$ws.Range("J94").Value = 5011.1665
$ws.Range("K94").Value = 223564.8
$ws.Range("L94").Value = 5011.1665
$ws.Range("M94").Value = -223113.8
$ws.Range("N94").Value = -5913.1665
$ws.Range("H105").Value = 2635.3076
$ws.Range("I105").Value = 1996
$ws.Range("J105").Value = 3183.2856
$ws.Range("K105").Value = 1996
$ws.Range("L105").Value = 3183.2856
$ws.Range("M105").Value = -249
$ws.Range("N105").Value = -6677.2856
$ws.Range("H132").Value = 6241.3076
$ws.Range("I132").Value = 6514.9
$ws.Range("K132").Value = 19544.7
$ws.Range("M132").Value = -17014.7
$ws.Range("H134").Value = 2865.3333
$ws.Range("I134").Value = 2838.4
$ws.Range("K134").Value = 8515.200000000001
$ws.Range("M134").Value = -5980.200000000001
$ws.Range("H136").Value = 1663.8
$ws.Range("J136").Value = 2164.5
$ws.Range("L136").Value = 6493.5
$ws.Range("N136").Value = -11593.5

# ---- Sheet: CUL ----
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H68").Value = 1000
$ws.Range("J68").Value = 1000
$ws.Range("L68").Value = 3000
$ws.Range("N68").Value = -4622
$ws.Range("H71").Value = 1000
$ws.Range("J71").Value = 1000
$ws.Range("L71").Value = 9000
$ws.Range("N71").Value = -17112
$ws.Range("H92").Value = 1950
$ws.Range("I92").Value = 900
$ws.Range("J92").Value = 3000
$ws.Range("K92").Value = 2700
$ws.Range("L92").Value = 9000
$ws.Range("M92").Value = -1452
$ws.Range("N92").Value = -11496
$ws.Range("H121").Value = 12828.667
$ws.Range("I121").Value = 16064.286
$ws.Range("J121").Value = 8298.799999999999
$ws.Range("K121").Value = 48192.858
$ws.Range("L121").Value = 24896.4
$ws.Range("M121").Value = -46882.858
$ws.Range("N121").Value = -27516.4
$ws.Range("H129").Value = 1251532
$ws.Range("I129").Value = 1370.5
$ws.Range("K129").Value = 4111.5
$ws.Range("M129").Value = 888.5

# ---- Sheet: GSM ----
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H80").Value = 3999
$ws.Range("I80").Value = 3999
$ws.Range("J80").Value = 0
$ws.Range("K80").Value = 3999
$ws.Range("L80").Value = 0
$ws.Range("M80").ClearContents()
$ws.Range("N80").Value = -3001
$ws.Range("H83").Value = 3999
$ws.Range("I83").Value = 3999
$ws.Range("J83").Value = 0
$ws.Range("K83").Value = 19995
$ws.Range("L83").Value = 0
$ws.Range("M83").ClearContents()
$ws.Range("N83").Value = -15003
$ws.Range("H102").Value = 1767.7273
$ws.Range("I102").Value = 1767.7273
$ws.Range("K102").Value = 1767.7273
$ws.Range("M102").Value = -145.7273
$ws.Range("H122").Value = 4389.364
$ws.Range("I122").Value = 4080.5
$ws.Range("K122").Value = 12241.5
$ws.Range("M122").Value = -9791.5
$ws.Range("H132").Value = 5542.7334
$ws.Range("I132").Value = 5748.2144
$ws.Range("K132").Value = 17244.6432
$ws.Range("M132").Value = -14714.6432
$ws.Range("H140").Value = 123495
$ws.Range("J140").Value = 123495
$ws.Range("L140").Value = 123495
$ws.Range("N140").Value = -133855

# ---- Sheet: LTW ----
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H33").Value = 1200
$ws.Range("I33").Value = 1200
$ws.Range("K33").Value = 1200
$ws.Range("M33").Value = -910
$ws.Range("H40").Value = 4449.6665
$ws.Range("I40").Value = 4151
$ws.Range("J40").Value = 4599
$ws.Range("K40").Value = 4151
$ws.Range("L40").Value = 4599
$ws.Range("M40").Value = -4015
$ws.Range("N40").Value = -4871
$ws.Range("H46").Value = 2101.4666
$ws.Range("I46").Value = 1728.3334
$ws.Range("J46").Value = 2661.1667
$ws.Range("K46").Value = 1728.3334
$ws.Range("L46").Value = 2661.1667
$ws.Range("M46").Value = -1540.3334
$ws.Range("N46").Value = -3037.1667
$ws.Range("H55").Value = 2027.3334
$ws.Range("I55").Value = 2425.3333
$ws.Range("J55").Value = 1828.3334
$ws.Range("K55").Value = 2425.3333
$ws.Range("L55").Value = 1828.3334
$ws.Range("M55").Value = -2252.3333
$ws.Range("N55").Value = -2174.3334
$ws.Range("H93").Value = 1625.5
$ws.Range("I93").Value = 1676
$ws.Range("K93").Value = 1676
$ws.Range("M93").Value = -428

# ---- Sheet: WVR ----
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H30").Value = 20000
$ws.Range("J30").Value = 20000
$ws.Range("L30").Value = 20000
$ws.Range("N30").Value = -20214
$ws.Range("H45").Value = 16010.286
$ws.Range("J45").Value = 16010.286
$ws.Range("L45").Value = 16010.286
$ws.Range("N45").Value = -16992.286
$ws.Range("H81").Value = 1001085.1
$ws.Range("I81").Value = 1205.4445
$ws.Range("K81").Value = 2410.889
$ws.Range("M81").Value = -1349.889
$ws.Range("H84").Value = 1001085.1
$ws.Range("I84").Value = 1205.4445
$ws.Range("K84").Value = 12054.445
$ws.Range("M84").Value = -6750.445
$ws.Range("H111").Value = 0
$ws.Range("J111").Value = 0
$ws.Range("L111").ClearContents()
$ws.Range("N111").Value = 0
$ws.Range("H132").Value = 1197.9166
$ws.Range("I132").Value = 1197.9166
$ws.Range("K132").Value = 3593.7498
$ws.Range("M132").Value = -1063.7498
